$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicated "Contact" / "No display for ContactDetail" row (row 10)
$meta.Rows.Item(10).Delete()

# Publisher value now populated
$meta.Range("B9").Value = "Alvearie Team"

# New Jurisdiction row (former second "Contact" row, now shifted up to row 10)
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet "Elements" (sheet2) ---
$elements = $wb.Worksheets.Item("Elements")

# Short / Definition for the root Extension slice summary row
$elements.Range("K2").Value = "Family Id"
$elements.Range("L2").Value = "The unique identifier for the subscriber (contract holder, employee) and the associated dependents"
